# Replace " дополнительной профессиональной программе " with a single space " ".
# This collapses the two runs that used to hold that text into the
# remaining run's content (" "), matching the target diff which merges
# the " дополнительной профессиональной" run and the " программе " run
# into a single run containing just " ".

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    " дополнительной профессиональной программе ",  # FindText
    $true,                                            # MatchCase
    $false,                                           # MatchWholeWord
    $false,                                           # MatchWildcards
    $false,                                           # MatchSoundsLike
    $false,                                           # MatchAllWordForms
    $true,                                             # Forward
    1,                                                 # Wrap (wdFindContinue)
    $false,                                            # Format
    " ",                                               # ReplaceWith
    2                                                  # Replace (wdReplaceAll)
)
